$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (1h % change) are plain text in this sheet: prices
# use dots as thousands separators (or are small decimals) and percentages
# carry a leading sign plus padding spaces. Several of the new Price values
# (e.g. "546.34", "12.00") look like ordinary numbers to Excel, so a bare
# assignment would silently convert them to numeric cells (and for values
# like "12.00"/"0.0790" would also drop the meaningful trailing zero).
# Prefixing those with a leading apostrophe -- exactly like typing '546.34
# into a cell -- forces Excel to keep them as text, matching the sheet.

$ws.Range("D2").Value = '61.505.09'
$ws.Range("E2").Value = '  -3.19%  '
$ws.Range("D3").Value = '2.999.93'
$ws.Range("E3").Value = '  -2.75%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''546.34'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = '''130.78'
$ws.Range("E6").Value = '  -6.14%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '2.996.04'
$ws.Range("E8").Value = '  -2.68%  '
$ws.Range("E9").Value = '  -1.80%  '
$ws.Range("E10").Value = '  -8.51%  '
$ws.Range("D11").Value = '''5.98'
$ws.Range("E11").Value = '  -6.83%  '
$ws.Range("D12").Value = '''0.444'
$ws.Range("E12").Value = '  -3.16%  '
$ws.Range("E13").Value = '  -3.15%  '
$ws.Range("D14").Value = '''33.98'
$ws.Range("D15").Value = '3.480.91'
$ws.Range("E15").Value = '  -2.87%  '
$ws.Range("D16").Value = '61.642.69'
$ws.Range("E16").Value = '  -2.95%  '
$ws.Range("E17").Value = '  -2.97%  '
$ws.Range("D18").Value = '2.997.36'
$ws.Range("E18").Value = '  -2.68%  '
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").Value = '''480.07'
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("D21").Value = '''13.16'
$ws.Range("E21").Value = '  -2.78%  '
$ws.Range("D22").Value = '''0.663'
$ws.Range("E23").Value = '  -1.90%  '
$ws.Range("D24").Value = '''80.65'
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("D25").Value = '''12.00'
$ws.Range("E25").Value = '  -2.00%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("D28").Value = '''7.64'
$ws.Range("E28").Value = '  -4.03%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").Value = '''25.47'
$ws.Range("E31").Value = '  -3.19%  '
$ws.Range("E32").Value = '  -4.07%  '
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").Value = '''5.52'
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").Value = '''54.89'
$ws.Range("E35").Value = '  -7.08%  '
$ws.Range("E36").Value = '  -3.09%  '
$ws.Range("D37").Value = '''446.87'
$ws.Range("E37").Value = '  -8.83%  '
$ws.Range("D38").Value = '3.124.16'
$ws.Range("E38").Value = '  -4.58%  '
$ws.Range("D39").Value = '''0.0790'
$ws.Range("E39").Value = '  -1.19%  '
$ws.Range("D40").Value = '''0.0381'
$ws.Range("E40").Value = '  -5.77%  '
$ws.Range("E41").Value = '  -2.05%  '
$ws.Range("D42").Value = '''8.06'
$ws.Range("E42").Value = '  -1.51%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '''2.35'
$ws.Range("E44").Value = '  -10.14%  '
$ws.Range("D45").Value = '''25.51'
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("E46").Value = '  -4.79%  '
$ws.Range("D47").Value = '''0.108'
$ws.Range("E47").Value = '  -1.73%  '
$ws.Range("D48").Value = '''1.94'
$ws.Range("E48").Value = '  -4.59%  '
$ws.Range("D49").Value = '''114.78'
$ws.Range("E49").Value = '  -7.67%  '
$ws.Range("E50").Value = '  +9.26%  '
$ws.Range("D51").Value = '0.0₃0483'
$ws.Range("E51").Value = '  -9.12%  '
